$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 14).Value = 282.76
$ws.Cells.Item(3, 14).Value = 226.46
$ws.Cells.Item(4, 14).Value = 173.64
$ws.Cells.Item(5, 14).Value = 133.93
$ws.Cells.Item(6, 14).Value = 127.31
$ws.Cells.Item(7, 14).Value = 125.27
$ws.Cells.Item(8, 14).Value = 109.04
$ws.Cells.Item(9, 14).Value = 73.04000000000001
$ws.Cells.Item(10, 14).Value = 63.85
$ws.Cells.Item(11, 14).Value = 48.88
$ws.Cells.Item(12, 14).Value = 47.55
$ws.Cells.Item(13, 14).Value = 40.59
$ws.Cells.Item(14, 14).Value = 39.04
$ws.Cells.Item(15, 14).Value = 37.86
$ws.Cells.Item(16, 14).Value = 36.96
$ws.Cells.Item(17, 14).Value = 35.49
$ws.Cells.Item(18, 14).Value = 33.33
$ws.Cells.Item(19, 14).Value = 23.78
$ws.Cells.Item(20, 14).Value = 19.71
$ws.Cells.Item(21, 14).Value = 16.45
$ws.Cells.Item(22, 14).Value = 15.62
$ws.Cells.Item(23, 14).Value = 13.12
$ws.Cells.Item(24, 14).Value = 12.82
$ws.Cells.Item(25, 14).Value = 11.11
$ws.Cells.Item(26, 14).Value = 10.37
$ws.Cells.Item(27, 14).Value = 9.56
$ws.Cells.Item(28, 14).Value = 5.93
$ws.Cells.Item(29, 14).Value = 5.56
$ws.Cells.Item(30, 14).Value = 5.56
$ws.Cells.Item(31, 14).Value = 3.63
$ws.Cells.Item(32, 14).Value = 3.33
$ws.Cells.Item(33, 14).Value = 3.33
$ws.Cells.Item(34, 14).Value = 3.33
$ws.Cells.Item(35, 14).Value = 3.33
$ws.Cells.Item(36, 14).Value = 3.33
$ws.Cells.Item(37, 14).Value = 2.78
$ws.Cells.Item(38, 14).Value = 0.97
$ws.Cells.Item(39, 14).Value = 0.6
$ws.Cells.Item(40, 14).Value = 0.45
$ws.Cells.Item(41, 14).Value = 0.3
$ws.Cells.Item(42, 14).Value = 0.3
$ws.Cells.Item(43, 14).Value = 0.3

$wb.Save()
